$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial value that was bumped by one day
# (45178 -> 45179) for every data row (rows 2 through 308).
$ws.Range("C2:C308").Value = 45179
